$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 32, pushing the existing "empty totals spacer" row
# (and the three summary rows below it) down by one.
$ws.Rows("32:32").Insert()

# Populate the newly-inserted row 32 with the next time-tracking entry
# (2014-03-01, 16:30 -> 18:00, i.e. 90 minutes / 1.5 hours).
$ws.Range("A32").Value = 2014
$ws.Range("B32").Value = 3
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 0.6875
$ws.Range("E32").Value = 0.75
$ws.Range("F32").Formula = "=(E32-D32)*24*60"
$ws.Range("G32").Formula = "=F32/60"

# Move the active selection to where the user clicked next.
$ws.Range("I37").Select()
